$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "MaxIsotopes" (row 8) into the new "IsotopeRange" parameter, which now
#    stores a min/max pair instead of a single number.
$ws.Range("A8").Value = "IsotopeRange"
$ws.Range("D8").Value = "The minimum and maximum number of isotopes to consider"

# 2. CorrelationMinimum (still row 5 at this point) "Update" column spelling fix:
#    Everytime -> Occasionally.
$ws.Range("C5").Value = "Occasionally"

# 3. Insert a new row above CorrelationMinimum (row 5) for the new
#    "MinAbsoluteChange" parameter, shifting CorrelationMinimum / IsotopicPercentage /
#    PPMThreshold / IsotopeRange / PlottingWindow / ProtonMass down by one.
$ws.Rows("5:5").Insert()

# 4. Populate the newly inserted row 5 with the new "MinAbsoluteChange" parameter.
$ws.Range("A5").Value = "MinAbsoluteChange"
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = "Occassionally"
$ws.Range("D5").Value = "An abundance (every peak is scaled to the largest peak) absolute change required to count a subsequent peak as an isotope. Default is 0.5."

# 5. NoiseFilter (row 3) default value changes from 5 to 1.
$ws.Range("B3").Value = 1

# 6. IsotopeRange is now row 9 after the insert; give it its "5,20" default,
#    stored as text so it keeps the "min,max" form.
$ws.Range("B9").Value = "5,20"
$ws.Range("B9").NumberFormat = "@"

# 7. PlottingWindow (now row 10) default value changes from 2 to 5.
$ws.Range("B10").Value = 5

# 8. Leave the selection on the new IsotopeRange default cell, as in the saved file.
$ws.Range("B9").Select()
